$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 1.0 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text fix (typo + trailing period) - all occurrences share the string
$ws.Range("B8").Value  = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B20").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B31").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B42").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B50").Value = "O usuário devidamente autenticado e na tela inicial do sistema."

# "histório" -> "histórico" typo fix
$ws.Range("B12").Value = "Chefe Verifica o histórico da tramitação da prestação de contas e clica para analisar a prestação de contas."
$ws.Range("B24").Value = "Chefe Verifica o histórico da tramitação da prestação de contas e clica para analisar a prestação de contas."
$ws.Range("B35").Value = "Chefe Verifica o histórico da tramitação da prestação de contas e clica para analisar a prestação de contas."
$ws.Range("B54").Value = "Chefe Verifica o histórico da tramitação da prestação de contas e clica para analisar a prestação de contas."

# Add trailing period to "SYSTEM Exibe a tela para prestação de contas"
$ws.Range("D12").Value = "SYSTEM Exibe a tela para prestação de contas."
$ws.Range("D24").Value = "SYSTEM Exibe a tela para prestação de contas."
$ws.Range("D35").Value = "SYSTEM Exibe a tela para prestação de contas."
$ws.Range("D54").Value = "SYSTEM Exibe a tela para prestação de contas."

# Add trailing period to "SYSTEM Exibe a tela de Detalhar Diárias"
$ws.Range("D26").Value = "SYSTEM Exibe a tela de Detalhar Diárias."

# Add semicolon before "Exibe mensagem de erro" in MSG203 text
$ws.Range("D56").Value = "SYSTEM Identifica que campos obrigatórios do parecer/análise não foram devidamente preenchidos; Exibe mensagem de erro (MSG203 - Campos obrigatórios) para o usuário."
